$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E value cells to be written as text (matching the source XML's inlineStr cells)
# rather than being auto-coerced to numbers by Excel's smart-entry parsing.
$valRange = $ws.Range("D2:E51")
$valRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.167.79'
$ws.Range("E2").Value = '  +3.94%  '
$ws.Range("D3").Value = '1.897.93'
$ws.Range("E3").Value = '  +4.16%  '
$ws.Range("D4").Value = '0.9956'
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '247.52'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").Value = '0.9968'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").Value = '0.4991'
$ws.Range("E7").Value = '  +1.21%  '
$ws.Range("D8").Value = '44.70'
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("D9").Value = '0.2959'
$ws.Range("E9").Value = '  +6.91%  '
$ws.Range("D10").Value = '0.06674'
$ws.Range("E10").Value = '  +4.51%  '
$ws.Range("D11").Value = '1.892.01'
$ws.Range("E11").Value = '  +3.82%  '
$ws.Range("D12").Value = '17.08'
$ws.Range("E12").Value = '  +2.58%  '
$ws.Range("D13").Value = '0.07204'
$ws.Range("E13").Value = '  +2.09%  '
$ws.Range("D14").Value = '0.6817'
$ws.Range("E14").Value = '  +6.00%  '
$ws.Range("D15").Value = '86.08'
$ws.Range("E15").Value = '  +2.39%  '
$ws.Range("D16").Value = '4.870'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("D17").Value = '30.131.87'
$ws.Range("E17").Value = '  +3.82%  '
$ws.Range("D18").Value = '0.000008062'
$ws.Range("E18").Value = '  +10.54%  '
$ws.Range("D19").Value = '0.9978'
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").Value = '12.98'
$ws.Range("E20").Value = '  +6.46%  '
$ws.Range("D21").Value = '2.134.76'
$ws.Range("E21").Value = '  +3.92%  '
$ws.Range("D22").Value = '0.9943'
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("D23").Value = '4.781'
$ws.Range("E23").Value = '  +5.15%  '
$ws.Range("D24").Value = '5.691'
$ws.Range("E24").Value = '  +6.02%  '
$ws.Range("D25").Value = '9.208'
$ws.Range("E25").Value = '  +4.21%  '
$ws.Range("D26").Value = '147.02'
$ws.Range("E26").Value = '  +2.42%  '
$ws.Range("D27").Value = '133.38'
$ws.Range("E27").Value = '  +2.90%  '
$ws.Range("D28").Value = '16.89'
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("D29").Value = '1.957'
$ws.Range("E29").Value = '  +4.03%  '
$ws.Range("D30").Value = '1.369'
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("D31").Value = '4.247'
$ws.Range("E31").Value = '  +2.94%  '
$ws.Range("D32").Value = '0.08763'
$ws.Range("E32").Value = '  +4.88%  '
$ws.Range("D33").Value = '3.964'
$ws.Range("E33").Value = '  +5.07%  '
$ws.Range("D34").Value = '0.05114'
$ws.Range("E34").Value = '  +3.23%  '
$ws.Range("D35").Value = '1.123'
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("D36").Value = '0.7109'
$ws.Range("E36").Value = '  +6.26%  '
$ws.Range("D37").Value = '2.659'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").Value = '2.808'
$ws.Range("E38").Value = '  +4.37%  '
$ws.Range("D39").Value = '2.240'
$ws.Range("E39").Value = '  -2.50%  '
$ws.Range("D40").Value = '0.9397'
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").Value = '0.01670'
$ws.Range("E41").Value = '  +5.42%  '
$ws.Range("D42").Value = '6.106'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").Value = '0.4239'
$ws.Range("E43").Value = '  +4.32%  '
$ws.Range("D44").Value = '0.9959'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").Value = '103.35'
$ws.Range("E45").Value = '  +2.30%  '
$ws.Range("D46").Value = '7.519'
$ws.Range("E46").Value = '  +4.91%  '
$ws.Range("D47").Value = '0.1266'
$ws.Range("E47").Value = '  +3.95%  '
$ws.Range("D48").Value = '0.05719'
$ws.Range("E48").Value = '  +3.22%  '
$ws.Range("D49").Value = '32.88'
$ws.Range("E49").Value = '  +3.80%  '
$ws.Range("D50").Value = '8.289'
$ws.Range("E50").Value = '  +2.66%  '
$ws.Range("D51").Value = '0.3759'
$ws.Range("E51").Value = '  +4.92%  '

# Clear the temporary text format so the cells end up with no explicit style,
# matching the original (unstyled) D/E cells.
$valRange.ClearFormats()
